{"js": "// Insert the contact-info line as a new centered paragraph directly after\n// the \"Dheeraj Chand\" name heading (first paragraph of the body).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the name heading paragraph (\"Dheeraj Chand\") \u2014 it's the first\n// paragraph in this resume document.\nlet nameParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Dheeraj Chand\") {\n    nameParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!nameParagraph) {\n  nameParagraph = paragraphs.items[0];\n}\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\n// Create the new paragraph right after the name heading. insertParagraph\n// clones the adjacent paragraph's run formatting (bold, 28pt), so clear the\n// inherited character formatting before writing the contact text \u2014 the\n// result should be a plain run, only the paragraph is centered.\nconst contactParagraph = nameParagraph.insertParagraph(\"\", \"After\");\ncontactParagraph.getRange().clear();\ncontactParagraph.alignment = Word.Alignment.centered;\ncontactParagraph.insertText(contactText, \"Start\");\n\nawait context.sync();\n", "ps1": "# Fix contact information missing from short resumes: add a centered\n# contact-info paragraph directly below the \"Dheeraj Chand\" name heading.\n$d = $word.ActiveDocument\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n# Word constants used below.\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n# Use Find/Replace with a paragraph-mark group character (^p) to split the\n# name heading into two paragraphs: the existing \"Dheeraj Chand\" run keeps\n# its original (bold, 28pt, centered) formatting, and the new contact-info\n# paragraph is created fresh \u2014 centered, but with no inherited direct\n# character formatting, matching a plain run typed after the heading.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Dheeraj Chand\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Dheeraj Chand^p\" + $contactText\n\n$find.Execute($find.Text, $False, $False, $False, $False, $False, $True, $wdFindContinue, $False, $find.Replacement.Text, $wdReplaceOne)\n"}
